# Apply the "Add 15 Product templates with correct industry content" edit:
# replace the AI/ML-flavoured sample copy with Product-Development copy
# across the three sheets of the Change Management Plan template.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Change Management Overview" ---------------------------------
$ws1 = $wb.Worksheets.Item("Change Management Overview")

$ws1.Range("A2").Value = "Product Development Implementation Project"
$ws1.Range("B6").Value = "Enterprise Product Development Implementation"
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new Product Development systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in Product Development technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for Product Development transformation"

# --- Sheet 2: "Change Impact Assessment" ------------------------------------
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

$ws2.Range("G4").Value = "Product Development automation"
$ws2.Range("G5").Value = "Product-powered insights"
$ws2.Range("G7").Value = "New Product interface"
$ws2.Range("G11").Value = "Product-enhanced CRM"
$ws2.Range("G12").Value = "Product-assisted support"
$ws2.Range("G13").Value = "Product-powered testing"
